$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ZY0606732"
$ws.Range("B3").Value = "ZY0606727"
$ws.Range("B4").Value = "ZY0605725"
$ws.Range("B5").Value = "ZY0602698"

$ws.Range("D5").Select()
